$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 2.88
$ws.Range("M5").Value = 1.11
$ws.Range("O5").Value = 1.54
$ws.Range("V5").Value = 1.1
$ws.Range("Y5").Value = 2.38
$ws.Range("Z5").Value = 1.53
$ws.Range("AN5").Value = 17
$ws.Range("AO5").Value = 51
$ws.Range("G6").Value = 2.7
$ws.Range("H6").Value = 2.88
$ws.Range("I6").Value = 3
$ws.Range("J6").Value = 3.6
$ws.Range("K6").Value = 1.8
$ws.Range("M6").Value = 1.14
$ws.Range("N6").Value = 5.5
$ws.Range("O6").Value = 1.62
$ws.Range("P6").Value = 2.2
$ws.Range("Q6").Value = 3.1
$ws.Range("R6").Value = 1.36
$ws.Range("S6").Value = 5.2
$ws.Range("U6").Value = 6.5
$ws.Range("V6").Value = 1.11
$ws.Range("W6").Value = 1.63
$ws.Range("Y6").Value = 2.38
$ws.Range("Z6").Value = 1.53
$ws.Range("AC6").Value = 12
$ws.Range("AD6").Value = 29
$ws.Range("AF6").Value = 51
$ws.Range("AJ6").Value = 101
$ws.Range("AN6").Value = 13
$ws.Range("H7").Value = 3.1
$ws.Range("W7").Value = 1.54
$ws.Range("AR7").Value = 1.9
$ws.Range("AS7").Value = 1.95
$ws.Range("G9").Value = 2.45
$ws.Range("H9").Value = 2.55
$ws.Range("I9").Value = 3.75
$ws.Range("L9").Value = 4.75
$ws.Range("M9").Value = 1.19
$ws.Range("N9").Value = 4
$ws.Range("P9").Value = 1.8
$ws.Range("Q9").Value = 4
$ws.Range("R9").Value = 1.25
$ws.Range("V9").Value = 1.05
$ws.Range("W9").Value = 1.8
$ws.Range("X9").Value = 1.95
$ws.Range("Y9").Value = 2.75
$ws.Range("Z9").Value = 1.4
$ws.Range("AA9").Value = 5
$ws.Range("AB9").Value = 9.5
$ws.Range("AE9").Value = 34
$ws.Range("AG9").Value = 4
$ws.Range("AI9").Value = 26
$ws.Range("AJ9").Value = 126
$ws.Range("AN9").Value = 17
$ws.Range("AQ9").Value = 67
$ws.Range("G14").Value = 1.48
$ws.Range("I14").Value = 6.25
$ws.Range("K14").Value = 2.4
$ws.Range("L14").Value = 7
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 10
$ws.Range("O14").Value = 1.36
$ws.Range("P14").Value = 3
$ws.Range("Q14").Value = 2.15
$ws.Range("R14").Value = 1.67
$ws.Range("S14").Value = 3.05
$ws.Range("T14").Value = 1.37
$ws.Range("V14").Value = 1.22
$ws.Range("Y14").Value = 2.38
$ws.Range("Z14").Value = 1.53
$ws.Range("AA14").Value = 5.5
$ws.Range("AB14").Value = 6
$ws.Range("AD14").Value = 9.5
$ws.Range("AF14").Value = 34
$ws.Range("AG14").Value = 8.5
$ws.Range("AL14").Value = 12
$ws.Range("AM14").Value = 29
$ws.Range("AQ14").Value = 51
$ws.Range("H15").Value = 3.9
$ws.Range("K15").Value = 2.2
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 1.3
$ws.Range("P15").Value = 3.4
$ws.Range("Q15").Value = 2
$ws.Range("R15").Value = 1.85
$ws.Range("S15").Value = 2.75
$ws.Range("T15").Value = 1.43
$ws.Range("U15").Value = 3.4
$ws.Range("V15").Value = 1.3
$ws.Range("W15").Value = 1.4
$ws.Range("X15").Value = 2.75
$ws.Range("Y15").Value = 2
$ws.Range("Z15").Value = 1.73
$ws.Range("AA15").Value = 6.5
$ws.Range("AE15").Value = 13
$ws.Range("AF15").Value = 29
$ws.Range("AG15").Value = 10
$ws.Range("AI15").Value = 19
$ws.Range("AJ15").Value = 51
$ws.Range("AK15").Value = 351
$ws.Range("AO15").Value = 51
$ws.Range("AP15").Value = 41
$ws.Range("G17").Value = 2.25
$ws.Range("I17").Value = 3.4
$ws.Range("J17").Value = 3
$ws.Range("L17").Value = 4
$ws.Range("N17").Value = 8
$ws.Range("Y17").Value = 1.83
$ws.Range("Z17").Value = 1.83
$ws.Range("AB17").Value = 10
$ws.Range("AD17").Value = 21
$ws.Range("AE17").Value = 21
$ws.Range("AL17").Value = 9.5
$ws.Range("AM17").Value = 17
$ws.Range("AO17").Value = 34
$ws.Range("AP17").Value = 29
$ws.Range("M28").Value = 1.07
$ws.Range("N28").Value = 7
$ws.Range("O28").Value = 1.41
$ws.Range("P28").Value = 2.62
$ws.Range("V28").Value = 1.13
$ws.Range("Z28").Value = 1.73
$ws.Range("AR28").Value = 1.85
$ws.Range("AS28").Value = 1.95
$ws.Range("G29").Value = 2.8
$ws.Range("I29").Value = 2.55
$ws.Range("J29").Value = 3.75
$ws.Range("L29").Value = 3.4
$ws.Range("M29").Value = 1.07
$ws.Range("N29").Value = 7
$ws.Range("O29").Value = 1.41
$ws.Range("P29").Value = 2.62
$ws.Range("Q29").Value = 2.4
$ws.Range("R29").Value = 1.53
$ws.Range("U29").Value = 4.5
$ws.Range("V29").Value = 1.15
$ws.Range("Z29").Value = 1.73
$ws.Range("AE29").Value = 29
$ws.Range("AO29").Value = 26
$ws.Range("AR29").Value = 1.83
$ws.Range("AS29").Value = 1.98
$ws.Range("G35").Value = 2.8
$ws.Range("I35").Value = 2.38
$ws.Range("L35").Value = 3.4
$ws.Range("Q35").Value = 2.6
$ws.Range("R35").Value = 1.48
$ws.Range("U35").Value = 5.5
$ws.Range("V35").Value = 1.14
$ws.Range("AA35").Value = 7
$ws.Range("AB35").Value = 13
$ws.Range("AC35").Value = 12
$ws.Range("AI35").Value = 21
$ws.Range("AK35").Value = 1250
$ws.Range("AL35").Value = 6
$ws.Range("AM35").Value = 10
$ws.Range("AO35").Value = 23
$ws.Range("AP35").Value = 23
$ws.Range("AR35").Value = 1.95
$ws.Range("AS35").Value = 1.85
$ws.Range("G37").Value = 1.33
$ws.Range("H37").Value = 4.75
$ws.Range("I37").Value = 9.5
$ws.Range("J37").Value = 1.83
$ws.Range("K37").Value = 2.38
$ws.Range("L37").Value = 9
$ws.Range("M37").Value = 1.06
$ws.Range("N37").Value = 10
$ws.Range("Q37").Value = 1.82
$ws.Range("R37").Value = 1.92
$ws.Range("Y37").Value = 2.25
$ws.Range("Z37").Value = 1.54
$ws.Range("AA37").Value = 6
$ws.Range("AB37").Value = 6
$ws.Range("AC37").Value = 9
$ws.Range("AD37").Value = 8
$ws.Range("AF37").Value = 34
$ws.Range("AH37").Value = 9.5
$ws.Range("AI37").Value = 23
$ws.Range("AJ37").Value = 81
$ws.Range("AL37").Value = 19
$ws.Range("AM37").Value = 41
$ws.Range("AN37").Value = 29
$ws.Range("AO37").Value = 126
$ws.Range("AP37").Value = 67
$ws.Range("AQ37").Value = 67
$ws.Range("I48").Value = 3.5
$ws.Range("K48").Value = 1.95
$ws.Range("Q48").Value = 2.4
$ws.Range("R48").Value = 1.53
$ws.Range("Z48").Value = 1.63
$ws.Range("AD48").Value = 19
$ws.Range("AE48").Value = 21
$ws.Range("AR48").Value = 1.77
$ws.Range("AS48").Value = 1.97
